# Update the division problems in the table to the new values.
$d = $word.ActiveDocument

$replacements = @(
    @{old="340÷5="; new="290÷5="},
    @{old="376÷5="; new="262÷6="},
    @{old="768÷8="; new="912÷6="},
    @{old="142÷3="; new="172÷3="},
    @{old="787÷7="; new="416÷7="},
    @{old="255÷3="; new="753÷5="},
    @{old="109÷2="; new="445÷6="},
    @{old="475÷6="; new="938÷3="},
    @{old="812÷2="; new="561÷5="},
    @{old="586÷5="; new="451÷7="},
    @{old="439÷5="; new="639÷3="},
    @{old="272÷4="; new="420÷2="},
    @{old="172÷8="; new="294÷6="},
    @{old="484÷7="; new="425÷4="},
    @{old="348÷6="; new="571÷9="},
    @{old="845÷5="; new="349÷8="},
    @{old="873÷3="; new="729÷6="},
    @{old="933÷8="; new="890÷4="},
    @{old="471÷3="; new="773÷2="},
    @{old="763÷3="; new="689÷6="},
    @{old="322÷8="; new="411÷5="},
    @{old="742÷9="; new="881÷5="},
    @{old="157÷6="; new="581÷2="},
    @{old="151÷8="; new="538÷5="},
    @{old="927÷9="; new="235÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
